$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9549892544746399
$ws.Range("B1").Value = 2.089877367019653
$ws.Range("C1").Value = 7.913414001464844
$ws.Range("D1").Value = 2.59632134437561
$ws.Range("E1").Value = 0.734183132648468
